# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Column G ("K") previously held total Strike# counts; replace with actual strikeout (K) counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 3
    3  = 5
    4  = 5
    5  = 9
    6  = 4
    7  = 10
    8  = 8
    9  = 3
    10 = 7
    11 = 4
    12 = 8
    13 = 12
    14 = 9
    15 = 8
    16 = 7
    17 = 6
    18 = 5
    19 = 10
    20 = 8
    21 = 8
    22 = 7
    23 = 7
    24 = 11
    25 = 7
    26 = 3
    27 = 8
    28 = 6
    29 = 10
    30 = 8
    31 = 6
    32 = 1
    33 = 6
    34 = 6
    35 = 6
    36 = 9
    37 = 7
    38 = 5
    39 = 2
    40 = 4
    41 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
